$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '30.229.87'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + '  -0.50%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'" + '2.063.16'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + '  +3.13%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'" + '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'" + '  +0.06%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'" + '325.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + '  +0.38%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'" + '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + '  +0.22%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'" + '0.5153'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'" + '  +1.11%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'" + '0.4299'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'" + '  +4.03%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'" + '0.08677'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'" + '  -0.54%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'" + '45.46'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + '  +5.47%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'" + '1.147'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + '  +1.27%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'" + '23.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'" + '  -2.20%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'" + '2.064.53'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + '  +3.47%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'" + '6.583'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + '  +0.34%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'" + '7.594'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'" + '  +1.96%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'" + '1.004'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + '  +0.12%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'" + '94.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'" + '  +0.21%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'" + '0.00001112'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'" + '  -0.08%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'" + '0.06599'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'" + '  +1.39%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'" + '18.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + '  -1.34%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'" + '  +0.22%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'" + '6.173'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + '  -0.32%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'" + '30.272.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + '  -0.57%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'" + '  +1.92%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'" + '2.272'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'" + '  +2.17%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'" + '2.304.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'" + '  +3.35%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'" + '21.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'" + '  -1.24%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'" + '161.80'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'" + '  -0.73%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'" + '2.479'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'" + '  +3.17%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'" + '129.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'" + '  -0.93%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'" + '1.162'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'" + '  +2.56%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'" + '  +0.86%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'" + '6.023'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'" + '  -0.60%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'" + '3.832'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'" + '  +0.19%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'" + '1.480'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'" + '  +10.49%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'" + '0.02540'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'" + '  +0.90%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'" + '9.499'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'" + '  +5.52%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'" + '5.372'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'" + '  -0.96%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'" + '0.06538'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'" + '  -0.85%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'" + '12.38'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + '  -0.67%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'" + '0.2217'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'" + '  +1.09%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'" + '0.6599'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + '  -0.29%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'" + '1.232'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'" + '  +0.13%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'" + '  +0.26%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'" + '13.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + '  +1.38%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'" + '0.6236'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'" + '  +1.35%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'" + '2.173'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'" + '  -0.52%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'" + '  -1.95%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'" + '1.227'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'" + '  -2.62%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'" + 'Aave'
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'" + 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'" + '80.83'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'" + '  +0.76%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'" + 'WEMIXTOKEN'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'" + 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'" + '1.169'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + '  +5.76%  '
$ws.Range("E51").Style = "Normal"
